# Apply crypto price/volume updates from the GitHub Actions data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.786.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.080.37'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.75'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.389'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0788'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.71%  '
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.386.56'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.97'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.762'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.086.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.704.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0823'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.140'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +10.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.94'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.41'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.119'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0626'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.73%  '
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.42'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0985'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.23%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.457.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("E47").Value = '  +3.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.30'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.83%  '

# Rows 48 and 49 swapped coins: InjectiveProtocol now ranks above FraxShare
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.43%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.43'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.31%  '
